# Auto-generated Excel COM-interop script to apply market-price refresh
# to the Carbuncle_Profits workbook (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 631.9655
$ws.Range("I80").Value = 221.92308
$ws.Range("K80").Value = 665.76924
$ws.Range("M80").Value = 332.23076
$ws.Range("H83").Value = 631.9655
$ws.Range("I83").Value = 221.92308
$ws.Range("K83").Value = 1997.30772
$ws.Range("M83").Value = 2994.69228
$ws.Range("H86").Value = 67592.11
$ws.Range("I86").Value = 34138.168
$ws.Range("J86").Value = 134500
$ws.Range("K86").Value = 34138.168
$ws.Range("L86").Value = 134500
$ws.Range("M86").Value = -33015.168
$ws.Range("N86").Value = -136746
$ws.Range("H89").Value = 67592.11
$ws.Range("I89").Value = 34138.168
$ws.Range("J89").Value = 134500
$ws.Range("K89").Value = 170690.84
$ws.Range("L89").Value = 672500
$ws.Range("M89").Value = -165074.84
$ws.Range("N89").Value = -683732
$ws.Range("H129").Value = 941.8
$ws.Range("I129").Value = 423
$ws.Range("K129").Value = 1269
$ws.Range("M129").Value = 3731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2408.4285
$ws.Range("I2").Value = 2408.4285
$ws.Range("K2").Value = 2408.4285
$ws.Range("M2").Value = -2295.4285
$ws.Range("H45").Value = 1523.3334
$ws.Range("J45").Value = 1586
$ws.Range("L45").Value = 1586
$ws.Range("N45").Value = -2340
$ws.Range("H74").Value = 1763.303
$ws.Range("I74").Value = 1287.4348
$ws.Range("J74").Value = 2857.8
$ws.Range("K74").Value = 1287.4348
$ws.Range("L74").Value = 2857.8
$ws.Range("M74").Value = -413.4348
$ws.Range("N74").Value = -4605.8
$ws.Range("H77").Value = 1763.303
$ws.Range("I77").Value = 1287.4348
$ws.Range("J77").Value = 2857.8
$ws.Range("K77").Value = 6437.174
$ws.Range("L77").Value = 14289
$ws.Range("M77").Value = -2069.174
$ws.Range("N77").Value = -23025
$ws.Range("H102").Value = 2854.7
$ws.Range("I102").Value = 2860.8333
$ws.Range("J102").Value = 2799.5
$ws.Range("K102").Value = 2860.8333
$ws.Range("L102").Value = 2799.5
$ws.Range("M102").Value = -1238.8333
$ws.Range("N102").Value = -6043.5
$ws.Range("H116").Value = 2408.4285
$ws.Range("I116").Value = 2408.4285
$ws.Range("K116").Value = 2408.4285
$ws.Range("M116").Value = -114.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2408.4285
$ws.Range("I3").Value = 2408.4285
$ws.Range("K3").Value = 2408.4285
$ws.Range("M3").Value = -2294.4285
$ws.Range("H113").Value = 1750
$ws.Range("I113").Value = 1750
$ws.Range("K113").Value = 1750
$ws.Range("M113").Value = 420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3726.5386
$ws.Range("I31").Value = 1561.0769
$ws.Range("J31").Value = 6613.8203
$ws.Range("K31").Value = 1561.0769
$ws.Range("L31").Value = 6613.8203
$ws.Range("M31").Value = -1266.0769
$ws.Range("N31").Value = -7203.8203
$ws.Range("H34").Value = 3726.5386
$ws.Range("I34").Value = 1561.0769
$ws.Range("J34").Value = 6613.8203
$ws.Range("K34").Value = 1561.0769
$ws.Range("L34").Value = 6613.8203
$ws.Range("M34").Value = -1359.0769
$ws.Range("N34").Value = -7017.8203
$ws.Range("H58").Value = 2397.7878
$ws.Range("I58").Value = 1889.5
$ws.Range("J58").Value = 2876.1765
$ws.Range("K58").Value = 1889.5
$ws.Range("L58").Value = 2876.1765
$ws.Range("M58").Value = -1686.5
$ws.Range("N58").Value = -3282.1765
$ws.Range("H86").Value = 19234510
$ws.Range("I86").Value = 25003110
$ws.Range("J86").Value = 5850
$ws.Range("K86").Value = 25003110
$ws.Range("L86").Value = 5850
$ws.Range("M86").Value = -25001987
$ws.Range("N86").Value = -8096
$ws.Range("H89").Value = 19234510
$ws.Range("I89").Value = 25003110
$ws.Range("J89").Value = 5850
$ws.Range("K89").Value = 125015550
$ws.Range("L89").Value = 29250
$ws.Range("M89").Value = -125009934
$ws.Range("N89").Value = -40482
$ws.Range("H136").Value = 2397.7878
$ws.Range("I136").Value = 1889.5
$ws.Range("J136").Value = 2876.1765
$ws.Range("K136").Value = 5668.5
$ws.Range("L136").Value = 8628.529500000001
$ws.Range("M136").Value = -3118.5
$ws.Range("N136").Value = -13728.5295

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 878707.5
$ws.Range("I5").Value = 758.5
$ws.Range("J5").Value = 1756656.5
$ws.Range("K5").Value = 2275.5
$ws.Range("L5").Value = 5269969.5
$ws.Range("M5").Value = -2163.5
$ws.Range("N5").Value = -5270193.5
$ws.Range("H128").Value = 125000
$ws.Range("I128").Value = 125000
$ws.Range("K128").Value = 375000
$ws.Range("M128").Value = -370020
$ws.Range("H131").Value = 843.2
$ws.Range("J131").Value = 870
$ws.Range("L131").Value = 2610
$ws.Range("N131").Value = -12690
$ws.Range("H133").Value = 2220
$ws.Range("I133").Value = 2513.3333
$ws.Range("K133").Value = 7539.999899999999
$ws.Range("M133").Value = -2479.999899999999
$ws.Range("H135").Value = 878707.5
$ws.Range("I135").Value = 758.5
$ws.Range("J135").Value = 1756656.5
$ws.Range("K135").Value = 6826.5
$ws.Range("L135").Value = 15809908.5
$ws.Range("M135").Value = -4291.5
$ws.Range("N135").Value = -15814978.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1912.8462
$ws.Range("I113").Value = 1722.1111
$ws.Range("J113").Value = 2342
$ws.Range("K113").Value = 1722.1111
$ws.Range("L113").Value = 2342
$ws.Range("M113").Value = 447.8888999999999
$ws.Range("N113").Value = -6682
$ws.Range("H126").Value = 2456.516
$ws.Range("I126").Value = 2336.6316
$ws.Range("J126").Value = 2646.3333
$ws.Range("K126").Value = 7009.8948
$ws.Range("L126").Value = 7938.999899999999
$ws.Range("M126").Value = -4539.8948
$ws.Range("N126").Value = -12878.9999
$ws.Range("H141").Value = 58227
$ws.Range("J141").Value = 58227
$ws.Range("L141").Value = 58227
$ws.Range("N141").Value = -68587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 5456.2
$ws.Range("I45").Value = 4500
$ws.Range("J45").Value = 5695.25
$ws.Range("K45").Value = 4500
$ws.Range("L45").Value = 5695.25
$ws.Range("M45").Value = -4009
$ws.Range("N45").Value = -6677.25
$ws.Range("H140").Value = 51161.145
$ws.Range("J140").Value = 52188
$ws.Range("L140").Value = 52188
$ws.Range("N140").Value = -62548
$ws.Range("H141").Value = 49512.5
$ws.Range("J141").Value = 49512.5
$ws.Range("L141").Value = 49512.5
$ws.Range("N141").Value = -59872.5

